$wb = $excel.ActiveWorkbook

$wsStats = $wb.Worksheets.Item("ugrad-009-01-stats-20")
$wsRaw = $wb.Worksheets.Item("raw")

# Reclassify 4 identifiers from "Context" to "Neither" for the first project row.
$wsRaw.Range("B2").Value = 47
$wsRaw.Range("C2").Value = 62

# Move selection on the stats sheet, then switch to/select the raw sheet.
$wsStats.Range("C3").Select()
$wsRaw.Activate()
$excel.Goto($wsRaw.Range("A2:C4"))

$wb.Save()
